$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.311.27"
$ws.Range("E2").Value = "  -0.65%  "
$ws.Range("D3").Value = "2.085.18"
$ws.Range("E3").Value = "  +3.45%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9987"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.42%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "328.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.03%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9980"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.33%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5197"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.63%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4324"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08843"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.74%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "46.20"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +6.50%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.160"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.13%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "24.61"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.28%  "
$ws.Range("D13").Value = "2.074.83"
$ws.Range("E13").Value = "  +3.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.689"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.22%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.701"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.71%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "95.22"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.70%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9986"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001121"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06612"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.27%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.83"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.9980"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.27%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.318"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.44%  "
$ws.Range("D23").Value = "30.341.36"
$ws.Range("E23").Value = "  -0.75%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.67%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.282"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.21%  "
$ws.Range("D26").Value = "2.319.74"
$ws.Range("E26").Value = "  +3.34%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.40"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.14%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.596"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.61%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "161.74"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "130.93"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.48%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.189"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.37%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.1067"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.23%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.642"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +20.87%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.253"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.818"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.34%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02582"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.24%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "9.801"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +7.62%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "12.70"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.12%  "
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06648"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.40%  "
$ws.Range("B40").Value = "InternetComputer(DFINITY)"
$ws.Range("C40").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.431"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.84%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.2251"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.39%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6810"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.245"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.18%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.9975"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.29%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.94"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.03%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6350"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.56%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.206"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.37%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.604"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.239"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -2.51%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.189"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +7.33%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "81.58"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.42%  "
